$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 5 (9d9de0ab... file): Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsZhCn.Range("D5").Value = "2016-01-28 05:57:31"
$wsZhCn.Range("G5").Value = "2016-01-28 05:58:40"

# de-de sheet, row 5 (9d9de0ab... file): Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
$wsDeDe.Range("D5").Value = "2016-01-28 05:58:05"
$wsDeDe.Range("G5").Value = "2016-01-28 05:58:58"
